$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 172
$ws.Cells.Item(172, 1).Value = 8
$ws.Cells.Item(172, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44595
$ws.Cells.Item(172, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(172, 5).Value = 4
$ws.Cells.Item(172, 6).Value = "Fruta"
$ws.Cells.Item(172, 7).Value = 100103
$ws.Cells.Item(172, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(172, 9).Value = 100103002
$ws.Cells.Item(172, 10).Value = "Ciruela"
$ws.Cells.Item(172, 11).Value = "Black Amber"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 20
$ws.Cells.Item(172, 14).Value = 240000
$ws.Cells.Item(172, 15).Value = 245000
$ws.Cells.Item(172, 16).Value = 242500
$ws.Cells.Item(172, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(172, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(172, 19).Value = 539
$ws.Cells.Item(172, 20).Value = 450

# New row 173
$ws.Cells.Item(173, 1).Value = 8
$ws.Cells.Item(173, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(173, 3).Value = "Coquimbo"
$ws.Cells.Item(173, 4).Value = 44595
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173, 5).Value = 4
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100103
$ws.Cells.Item(173, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(173, 9).Value = 100103002
$ws.Cells.Item(173, 10).Value = "Ciruela"
$ws.Cells.Item(173, 11).Value = "Black Amber"
$ws.Cells.Item(173, 12).Value = "Segunda"
$ws.Cells.Item(173, 13).Value = 24
$ws.Cells.Item(173, 14).Value = 195000
$ws.Cells.Item(173, 15).Value = 200000
$ws.Cells.Item(173, 16).Value = 197500
$ws.Cells.Item(173, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(173, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value = 439
$ws.Cells.Item(173, 20).Value = 450
